$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text updates (issue number 49 -> 50, week-of dates refreshed)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# ---------------------------------------------------------------------------
# 2) Cells whose TYPE changes (number <-> text placeholder) need their style
#    copied from a matching, untouched donor cell before the new value is
#    written (PasteSpecial formats keeps the exact existing style index
#    instead of minting a new one). Donor cells used (never modified below):
#      C22 = text "0"      (s=14)
#      E22 = text "***.*"  (s=14)
#      F27 = number, s=15
#      K27 = number, s=16
# ---------------------------------------------------------------------------

# --- L14: "***.*" (text) -> 0 (number, style like K14/M14/N14) ---
$ws.Range("K27").Copy()
$ws.Range("L14").PasteSpecial(-4122)
$ws.Range("L14").Value = 0

# --- C15: 1 (number) -> "0" (text) ---
$ws.Range("C22").Copy()
$ws.Range("C15").PasteSpecial(-4104)
$ws.Range("C22").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# --- Row 16: columns shift by one (C/D/E change shape) ---
# C16: "0" (text) -> 1 (number)
$ws.Range("F27").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 1
# D16: 2 (number) -> "0" (text)
$ws.Range("C22").Copy()
$ws.Range("D16").PasteSpecial(-4104)
$ws.Range("C22").Copy()
$ws.Range("D16").PasteSpecial(-4122)
# E16: -100 (number) -> "***.*" (text)
$ws.Range("E22").Copy()
$ws.Range("E16").PasteSpecial(-4104)
$ws.Range("E22").Copy()
$ws.Range("E16").PasteSpecial(-4122)
# F16: 1 -> 2 (plain value update, style unchanged)
$ws.Range("F16").Value = 2

# --- Row 17 ---
# C17: 1 (number) -> "0" (text)
$ws.Range("C22").Copy()
$ws.Range("C17").PasteSpecial(-4104)
$ws.Range("C22").Copy()
$ws.Range("C17").PasteSpecial(-4122)

# --- Row 18 ---
# C18: 4 (number) -> "0" (text)
$ws.Range("C22").Copy()
$ws.Range("C18").PasteSpecial(-4104)
$ws.Range("C22").Copy()
$ws.Range("C18").PasteSpecial(-4122)
# D18: "0" (text) -> 1 (number)
$ws.Range("F27").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 1
# E18: "***.*" (text) -> -100 (number)
$ws.Range("K27").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = -100

# --- C26: 1 (number) -> "0" (text) ---
$ws.Range("C22").Copy()
$ws.Range("C26").PasteSpecial(-4104)
$ws.Range("C22").Copy()
$ws.Range("C26").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Plain value-only updates (style / type unchanged)
# ---------------------------------------------------------------------------

# Row 16 (remaining recomputed figures)
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 20
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -9.090909090909
$ws.Range("N16").Value = -68.253968253968

# Row 17 (remaining recomputed figures)
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 4
$ws.Range("H17").Value = -33.333333333333
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 70
$ws.Range("M17").Value = 77.083333333333
$ws.Range("N17").Value = -19.047619047619

# Row 18 (remaining recomputed figures)
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 30.952380952381
$ws.Range("L18").Value = 129.166666666667
$ws.Range("M18").Value = -47.115384615384
$ws.Range("N18").Value = -83.180428134556

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -17.391304347826
$ws.Range("I19").Value = 275
$ws.Range("J19").Value = 263
$ws.Range("K19").Value = 4.562737642585
$ws.Range("L19").Value = 41.752577319587
$ws.Range("M19").Value = 90.972222222222
$ws.Range("N19").Value = 17.521367521367

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 78
$ws.Range("J20").Value = 119
$ws.Range("K20").Value = -34.453781512605
$ws.Range("L20").Value = 32.203389830508
$ws.Range("M20").Value = 110.810810810811
$ws.Range("N20").Value = -88.793103448275

# Row 21 (TOTAL)
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 38
$ws.Range("G21").Value = 48
$ws.Range("H21").Value = -20.833333333333
$ws.Range("I21").Value = 516
$ws.Range("J21").Value = 498
$ws.Range("K21").Value = 3.614457831325
$ws.Range("L21").Value = 60.747663551401
$ws.Range("M21").Value = 43.732590529247
$ws.Range("N21").Value = -63.916083916083

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 52
$ws.Range("G24").Value = 44
$ws.Range("H24").Value = 18.181818181818
$ws.Range("I24").Value = 492
$ws.Range("J24").Value = 480
$ws.Range("K24").Value = 2.5
$ws.Range("L24").Value = 78.909090909090
$ws.Range("M24").Value = -6.285714285714

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 33.333333333333
$ws.Range("I25").Value = 187
$ws.Range("J25").Value = 177
$ws.Range("K25").Value = 5.649717514124
$ws.Range("L25").Value = 33.571428571428
$ws.Range("M25").Value = -12.206572769953

# Row 28 (Shooting Vic.)
$ws.Range("L28").Value = 0

# Row 29 (Shooting Inc.)
$ws.Range("L29").Value = 0
